$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accessories")

$cell = $ws.Range("A6")
$cell.Value = "Cases & Protection"
$cell.WrapText = $true

$ws.Activate()
$cell.Select() | Out-Null
